$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rich-text edits (shared strings) ---

# A8: "Volume 29   Number  47" -> "...48"
$a8 = $ws.Range("A8")
$a8.Characters(21,2).Text = "48"
$a8.Characters(1,21).Font.Name = "Andale WT"
$a8.Characters(1,21).Font.Size = 10
$a8.Characters(22,1).Font.Name = "Andale WT"
$a8.Characters(22,1).Font.Size = 10

# C9: "Report Covering the Week  11/21/2022  Through  11/27/2022"
#  -> "Report Covering the Week  11/28/2022  Through  12/4/2022"
$c9 = $ws.Range("C9")
$c9.Characters(27,10).Text = "11/28/2022"
$c9.Characters(48,10).Text = "12/4/2022"
$c9full = $c9.Text
$c9len = $c9full.Length
$c9.Characters(1, $c9len - 1).Font.Name = "Andale WT"
$c9.Characters(1, $c9len - 1).Font.Size = 10
$c9.Characters($c9len, 1).Font.Name = "Andale WT"
$c9.Characters($c9len, 1).Font.Size = 10

# --- Cell C15 / C30: shared-string placeholder "0" -> real numeric value ---
# Convert these from text "0" to numbers, matching sibling-column number format (#,##0).
$ws.Range("C15").Value = 4
$ws.Range("C15").NumberFormat = "#,##0"

$ws.Range("C30").Value = 1
$ws.Range("C30").NumberFormat = "#,##0"

# --- Remaining numeric cell updates (rows 14-30) ---
$ws.Range("G14").Value = 2
$ws.Range("H14").Value = 50
$ws.Range("J14").Value = 87
$ws.Range("K14").Value = -18.390804597701
$ws.Range("N14").Value = -83.52668213457
$ws.Range("D15").Value = 7
$ws.Range("E15").Value = -42.857142857142
$ws.Range("F15").Value = 14
$ws.Range("G15").Value = 23
$ws.Range("H15").Value = -39.130434782608
$ws.Range("I15").Value = 234
$ws.Range("J15").Value = 210
$ws.Range("K15").Value = 11.428571428571
$ws.Range("L15").Value = 9.859154929577
$ws.Range("M15").Value = 6.849315068493
$ws.Range("N15").Value = -58.802816901408
$ws.Range("C16").Value = 48
$ws.Range("D16").Value = 54
$ws.Range("E16").Value = -11.111111111111
$ws.Range("F16").Value = 189
$ws.Range("G16").Value = 186
$ws.Range("H16").Value = 1.612903225806
$ws.Range("I16").Value = 2384
$ws.Range("J16").Value = 1946
$ws.Range("K16").Value = 22.507708119218
$ws.Range("L16").Value = 19.498746867167
$ws.Range("M16").Value = -30.394160583941
$ws.Range("N16").Value = -84.723824170191
$ws.Range("C17").Value = 60
$ws.Range("D17").Value = 67
$ws.Range("E17").Value = -10.447761194029
$ws.Range("F17").Value = 273
$ws.Range("G17").Value = 286
$ws.Range("H17").Value = -4.545454545454
$ws.Range("I17").Value = 3834
$ws.Range("J17").Value = 3344
$ws.Range("K17").Value = 14.653110047846
$ws.Range("L17").Value = 20.641913152926
$ws.Range("M17").Value = 24.359390204346
$ws.Range("N17").Value = -51.566447700859
$ws.Range("C18").Value = 37
$ws.Range("D18").Value = 47
$ws.Range("E18").Value = -21.27659574468
$ws.Range("F18").Value = 149
$ws.Range("G18").Value = 192
$ws.Range("H18").Value = -22.395833333333
$ws.Range("I18").Value = 2198
$ws.Range("J18").Value = 1941
$ws.Range("K18").Value = 13.240597630087
$ws.Range("L18").Value = -12.5
$ws.Range("M18").Value = -27.025232403718
$ws.Range("N18").Value = -80.418708240534
$ws.Range("C19").Value = 110
$ws.Range("D19").Value = 128
$ws.Range("E19").Value = -14.0625
$ws.Range("F19").Value = 421
$ws.Range("G19").Value = 465
$ws.Range("H19").Value = -9.462365591397
$ws.Range("I19").Value = 5548
$ws.Range("J19").Value = 4427
$ws.Range("K19").Value = 25.321888412017
$ws.Range("L19").Value = 29.023255813953
$ws.Range("M19").Value = 37.462834489593
$ws.Range("N19").Value = -14.091049860638
$ws.Range("C20").Value = 26
$ws.Range("D20").Value = 24
$ws.Range("E20").Value = 8.333333333333
$ws.Range("F20").Value = 138
$ws.Range("G20").Value = 121
$ws.Range("H20").Value = 14.049586776859
$ws.Range("I20").Value = 1721
$ws.Range("J20").Value = 1426
$ws.Range("K20").Value = 20.687237026648
$ws.Range("L20").Value = 34.348165495706
$ws.Range("M20").Value = 29.495861550037
$ws.Range("N20").Value = -80.643347205038
$ws.Range("C21").Value = 285
$ws.Range("D21").Value = 328
$ws.Range("E21").Value = -13.109756097561
$ws.Range("F21").Value = 1187
$ws.Range("G21").Value = 1275
$ws.Range("H21").Value = -6.901960784313
$ws.Range("I21").Value = 15990
$ws.Range("J21").Value = 13381
$ws.Range("K21").Value = 19.497795381511
$ws.Range("L21").Value = 17.746686303387
$ws.Range("M21").Value = 4.9970451113
$ws.Range("N21").Value = -68.705352774244
$ws.Range("D22").Value = 9
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 32
$ws.Range("G22").Value = 34
$ws.Range("H22").Value = -5.882352941176
$ws.Range("I22").Value = 326
$ws.Range("J22").Value = 261
$ws.Range("K22").Value = 24.904214559387
$ws.Range("L22").Value = 10.135135135135
$ws.Range("M22").Value = -19.901719901719
$ws.Range("C23").Value = 19
$ws.Range("D23").Value = 34
$ws.Range("E23").Value = -44.117647058823
$ws.Range("F23").Value = 98
$ws.Range("G23").Value = 125
$ws.Range("H23").Value = -21.6
$ws.Range("I23").Value = 1396
$ws.Range("J23").Value = 1371
$ws.Range("K23").Value = 1.823486506199
$ws.Range("L23").Value = 8.807482462977
$ws.Range("M23").Value = 26.334841628959
$ws.Range("C24").Value = 274
$ws.Range("D24").Value = 260
$ws.Range("E24").Value = 5.384615384615
$ws.Range("F24").Value = 964
$ws.Range("G24").Value = 1001
$ws.Range("H24").Value = -3.696303696303
$ws.Range("I24").Value = 12498
$ws.Range("J24").Value = 9794
$ws.Range("K24").Value = 27.608740044925
$ws.Range("L24").Value = 25.105105105105
$ws.Range("M24").Value = 27.595712098009
$ws.Range("C25").Value = 77
$ws.Range("D25").Value = 100
$ws.Range("E25").Value = -23
$ws.Range("F25").Value = 401
$ws.Range("G25").Value = 424
$ws.Range("H25").Value = -5.424528301886
$ws.Range("I25").Value = 5485
$ws.Range("J25").Value = 4464
$ws.Range("K25").Value = 22.871863799283
$ws.Range("L25").Value = 28.665259207131
$ws.Range("M25").Value = -25.384301455584
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = -54.545454545454
$ws.Range("F26").Value = 26
$ws.Range("G26").Value = 33
$ws.Range("H26").Value = -21.212121212121
$ws.Range("I26").Value = 356
$ws.Range("J26").Value = 359
$ws.Range("K26").Value = -0.8356545961
$ws.Range("L26").Value = 7.878787878787
$ws.Range("C27").Value = 4
$ws.Range("D27").Value = 20
$ws.Range("E27").Value = -80
$ws.Range("F27").Value = 36
$ws.Range("G27").Value = 62
$ws.Range("H27").Value = -41.935483870967
$ws.Range("I27").Value = 573
$ws.Range("J27").Value = 646
$ws.Range("K27").Value = -11.300309597523
$ws.Range("L27").Value = 3.057553956834
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 7
$ws.Range("E28").Value = -85.714285714285
$ws.Range("F28").Value = 15
$ws.Range("H28").Value = -34.782608695652
$ws.Range("I28").Value = 319
$ws.Range("J28").Value = 390
$ws.Range("K28").Value = -18.205128205128
$ws.Range("L28").Value = -32.415254237288
$ws.Range("M28").Value = -33.123689727463
$ws.Range("N28").Value = -81.624423963133
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 7
$ws.Range("E29").Value = -85.714285714285
$ws.Range("F29").Value = 14
$ws.Range("G29").Value = 22
$ws.Range("H29").Value = -36.363636363636
$ws.Range("I29").Value = 268
$ws.Range("J29").Value = 319
$ws.Range("K29").Value = -15.987460815047
$ws.Range("L29").Value = -31.282051282051
$ws.Range("M29").Value = -30.927835051546
$ws.Range("N29").Value = -82.820512820512
$ws.Range("F30").Value = 9
$ws.Range("H30").Value = 350
$ws.Range("I30").Value = 84
$ws.Range("K30").Value = 50
$ws.Range("L30").Value = 115.384615384615
